$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.29499999999999
$ws.Range("D6").Value = -7.894800000000004
$ws.Range("D7").Value = -7.532799999999996
$ws.Range("C8").Value = -12.49769999999999
$ws.Range("D8").Value = -8.498699999999989
$ws.Range("E11").Value = 13.7216
$ws.Range("B12").Value = 6.058199999999997
$ws.Range("C12").Value = -14.54590000000002
$ws.Range("C14").Value = -12.361
$ws.Range("E14").Value = 14.02020000000001
$ws.Range("D19").Value = -7.683699999999996
$ws.Range("E19").Value = 14.1191
$ws.Range("D21").Value = -7.445400000000003
$ws.Range("E21").Value = 14.1515
$ws.Range("C22").Value = -11.04309999999999
$ws.Range("D24").Value = -8.398399999999992
